$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 89
$ws.Range("A89").Value = '05/01/2026 13:25:36'
$ws.Range("B89").Value = '05/01 13:12'
$ws.Range("C89").Value = 'Metrópoles'
$ws.Range("D89").Value = 'Projeto sugere título de Cidadão Honorário de Brasília a Nunes Marques'
$ws.Range("E89").Value = 'https://www.metropoles.com/colunas/grande-angular/projeto-sugere-titulo-de-cidadao-honorario-de-brasilia-a-nunes-marques'
$ws.Range("F89").Value = 'stf'
$ws.Range("G89").Value = 'Kassio Nunes Marques nasceu em Teresina (PI). O ministro está no STF desde 2020 e, em 2026, assumirá também a função de presidente do TSE'

# Row 90
$ws.Range("A90").Value = '05/01/2026 13:25:37'
$ws.Range("B90").Value = '05/01 13:05'
$ws.Range("C90").Value = 'Metrópoles'
$ws.Range("D90").Value = 'Em despacho, ministro do TCU prevê reversão da liquidação do Master'
$ws.Range("E90").Value = 'https://www.metropoles.com/colunas/tacio-lorran/em-despacho-ministro-do-tcu-preve-reversao-da-liquidacao-do-master'
$ws.Range("F90").Value = 'tcu'
$ws.Range("G90").Value = 'Ministro do TCU, Jhonatan de Jesus determinou que TCU inspeciona documentos do caso Master em posse do Ba'

# Row 91
$ws.Range("A91").Value = '05/01/2026 13:25:38'
$ws.Range("B91").Value = '05/01 13:00'
$ws.Range("C91").Value = 'Folha de S.Paulo - Poder - Principal'
$ws.Range("D91").Value = 'Lula planeja evento no 8/1 para anunciar veto a projeto que reduz pena de Bolsonaro'
$ws.Range("E91").Value = 'https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/poder/2026/01/lula-planeja-evento-no-81-para-anunciar-veto-a-projeto-que-reduz-pena-de-bolsonaro.shtml'
$ws.Range("F91").Value = 'lula'
$ws.Range("G91").Value = 'ade de anúncio do veto do presidente &lt;a href="https://www1.folha.uol.com.br/folha-topicos/&lt;b&gt;lula&lt;/b&gt;/"&gt;Lula&lt;/a&gt; (PT) à redução de penas dos condenados por atos golpistas, entre eles &lt;a href='
